# add guild data module
# Adds a new row (row 11) to the "Property" sheet describing a "GuilID"
# field, mirroring the existing rows' shape/formatting, and moves the
# sheet's active selection to E19 (as recorded by the original author).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$ws.Cells.Item(11, 1).Value  = "GuilID"
$ws.Cells.Item(11, 2).Value  = "object"
$ws.Cells.Item(11, 3).Value  = $true
$ws.Cells.Item(11, 4).Value  = $true
$ws.Cells.Item(11, 5).Value  = $true
$ws.Cells.Item(11, 6).Value  = $true
$ws.Cells.Item(11, 7).Value  = 0
$ws.Cells.Item(11, 8).Value  = 0
$ws.Cells.Item(11, 9).Value  = "Friend"
$ws.Cells.Item(11, 10).Value = "工会ID"

# Columns A, B, I, J use the text ("@") number format on every other row
# in this sheet - match that so the new row's style lines up too.
$ws.Cells.Item(11, 1).NumberFormat  = "@"
$ws.Cells.Item(11, 2).NumberFormat  = "@"
$ws.Cells.Item(11, 9).NumberFormat  = "@"
$ws.Cells.Item(11, 10).NumberFormat = "@"

# Restore the active selection recorded in the saved workbook.
$ws.Range("E19").Select()
